$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 3.230985683306322
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 0.1575252929769615
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 13.71653804550039

# Row 3
$ws.Range("B3").Value = 0.003994804209775715
$ws.Range("C3").Value = 0.04240448674262143
$ws.Range("D3").Value = 3.900430680208489
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 4.443609181331619

# Row 4
$ws.Range("B4").Value = 3.230985683306322
$ws.Range("C4").Value = 1.667794583268128
$ws.Range("D4").Value = 26.21740644021617
$ws.Range("E4").Value = 0.496779210170732
$ws.Range("G4").Value = 31.61296591696135
